$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value (avoids Excel auto-converting
# numeric-looking strings like "3.39" or "1.00" into actual numbers).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '69.608.23'
Set-TextValue $ws.Range("E2") '  +1.92%  '

Set-TextValue $ws.Range("D3") '3.936.08'
Set-TextValue $ws.Range("E3") '  +0.56%  '

Set-TextValue $ws.Range("E4") '  +0.05%  '

Set-TextValue $ws.Range("D5") '531.29'
Set-TextValue $ws.Range("E5") '  +8.74%  '

Set-TextValue $ws.Range("D6") '145.25'
Set-TextValue $ws.Range("E6") '  -1.06%  '

Set-TextValue $ws.Range("D7") '0.620'
Set-TextValue $ws.Range("E7") '  -0.31%  '

Set-TextValue $ws.Range("E8") '  +0.00%  '

Set-TextValue $ws.Range("D9") '0.730'
Set-TextValue $ws.Range("E9") '  +0.26%  '

Set-TextValue $ws.Range("E10") '  +4.34%  '

Set-TextValue $ws.Range("E11") '  -1.31%  '

Set-TextValue $ws.Range("D12") '42.72'
Set-TextValue $ws.Range("E12") '  -0.79%  '

Set-TextValue $ws.Range("D13") '10.43'
Set-TextValue $ws.Range("E13") '  -3.67%  '

Set-TextValue $ws.Range("D14") '4.565.34'
Set-TextValue $ws.Range("E14") '  +0.59%  '

Set-TextValue $ws.Range("D15") '3.944.72'
Set-TextValue $ws.Range("E15") '  +1.02%  '

Set-TextValue $ws.Range("D16") '14.05'
Set-TextValue $ws.Range("E16") '  -1.03%  '

Set-TextValue $ws.Range("E17") '  -0.23%  '

Set-TextValue $ws.Range("E18") '  +6.84%  '

Set-TextValue $ws.Range("D19") '19.89'
Set-TextValue $ws.Range("E19") '  +0.04%  '

Set-TextValue $ws.Range("D20") '69.490.06'
Set-TextValue $ws.Range("E20") '  +1.63%  '

Set-TextValue $ws.Range("D21") '432.88'
Set-TextValue $ws.Range("E21") '  -0.29%  '

Set-TextValue $ws.Range("D22") '3.39'
Set-TextValue $ws.Range("E22") '  -3.86%  '

Set-TextValue $ws.Range("D23") '14.58'
Set-TextValue $ws.Range("E23") '  -2.26%  '

Set-TextValue $ws.Range("D24") '88.67'
Set-TextValue $ws.Range("E24") '  +1.03%  '

Set-TextValue $ws.Range("D25") '4.11'
Set-TextValue $ws.Range("E25") '  +14.05%  '

Set-TextValue $ws.Range("D26") '11.92'
Set-TextValue $ws.Range("E26") '  +3.96%  '

Set-TextValue $ws.Range("D27") '10.84'
Set-TextValue $ws.Range("E27") '  -3.71%  '

Set-TextValue $ws.Range("D28") '36.69'
Set-TextValue $ws.Range("E28") '  -3.95%  '

Set-TextValue $ws.Range("D29") '702.15'
Set-TextValue $ws.Range("E29") '  -3.28%  '

Set-TextValue $ws.Range("E30") '  -2.91%  '

Set-TextValue $ws.Range("E31") '  -1.63%  '

Set-TextValue $ws.Range("D32") '2.87'
Set-TextValue $ws.Range("E32") '  -1.70%  '

Set-TextValue $ws.Range("D33") '69.05'
Set-TextValue $ws.Range("E33") '  +13.86%  '

Set-TextValue $ws.Range("D34") '0.455'
Set-TextValue $ws.Range("E34") '  +12.34%  '

Set-TextValue $ws.Range("D35") '6.13'
Set-TextValue $ws.Range("E35") '  -2.14%  '

Set-TextValue $ws.Range("D36") '0.0₃0868'
Set-TextValue $ws.Range("E36") '  -0.35%  '

Set-TextValue $ws.Range("D37") '40.56'
Set-TextValue $ws.Range("E37") '  -2.65%  '

Set-TextValue $ws.Range("D38") '0.149'
Set-TextValue $ws.Range("E38") '  +1.02%  '

Set-TextValue $ws.Range("E39") '  -0.01%  '

Set-TextValue $ws.Range("E40") '  -0.05%  '

Set-TextValue $ws.Range("D41") '0.0483'
Set-TextValue $ws.Range("E41") '  +0.45%  '

Set-TextValue $ws.Range("D42") '2.83'
Set-TextValue $ws.Range("E42") '  -4.08%  '

Set-TextValue $ws.Range("D43") '3.10'
Set-TextValue $ws.Range("E43") '  +6.51%  '

Set-TextValue $ws.Range("D44") '3.00'
Set-TextValue $ws.Range("E44") '  -4.86%  '

Set-TextValue $ws.Range("D45") '3.23'
Set-TextValue $ws.Range("E45") '  +14.69%  '

Set-TextValue $ws.Range("E46") '  +2.39%  '

Set-TextValue $ws.Range("D47") '0.142'
Set-TextValue $ws.Range("E47") '  +1.09%  '

Set-TextValue $ws.Range("D48") '0.0₆0358'
Set-TextValue $ws.Range("E48") '  +2.88%  '

Set-TextValue $ws.Range("E49") '  -2.04%  '

Set-TextValue $ws.Range("E50") '  -1.66%  '

Set-TextValue $ws.Range("D51") '144.74'
Set-TextValue $ws.Range("E51") '  -0.11%  '
